$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells are formatted as text so values like '34.667.21' or '1.00' are not reinterpreted as numbers/dates

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value2 = '34.667.21'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value2 = '  +0.92%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value2 = '1.805.78'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value2 = '  +0.13%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value2 = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value2 = '225.07'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value2 = '  -1.04%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value2 = '  +0.50%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value2 = '  -0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value2 = '39.35'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value2 = '  +8.75%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value2 = '0.290'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value2 = '  -2.23%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value2 = '  -3.07%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value2 = '0.100'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value2 = '  +3.65%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value2 = '2.068.96'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value2 = '  +0.00%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value2 = '1.808.67'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value2 = '  -0.34%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value2 = '10.93'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value2 = '  -3.06%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value2 = '  -1.29%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value2 = '34.672.30'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value2 = '  +0.78%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value2 = '4.38'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value2 = '  -1.28%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value2 = '  -2.68%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value2 = '240.86'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value2 = '  -1.57%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value2 = '0.0₃0769'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value2 = '  -1.90%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value2 = '11.11'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value2 = '  -2.83%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value2 = '1.00'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value2 = '  +0.10%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value2 = '  -1.91%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value2 = '  -2.61%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value2 = '171.72'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value2 = '  +0.51%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value2 = '7.67'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value2 = '  -5.57%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value2 = '17.44'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value2 = '  +0.23%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value2 = '  -0.14%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value2 = '  +0.01%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value2 = '  -1.54%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value2 = '3.76'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value2 = '  -1.30%  '
$ws.Range('B32').NumberFormat = '@'
$ws.Range('B32').Value2 = 'Hedera'
$ws.Range('C32').NumberFormat = '@'
$ws.Range('C32').Value2 = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value2 = '0.0515'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value2 = '  -1.33%  '
$ws.Range('B33').NumberFormat = '@'
$ws.Range('B33').Value2 = 'InternetComputer(DFINITY)'
$ws.Range('C33').NumberFormat = '@'
$ws.Range('C33').Value2 = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value2 = '3.85'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value2 = '  -3.18%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value2 = '1.81'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value2 = '  +1.13%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value2 = '0.643'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value2 = '  -2.87%  '
$ws.Range('B36').NumberFormat = '@'
$ws.Range('B36').Value2 = 'Maker'
$ws.Range('C36').NumberFormat = '@'
$ws.Range('C36').Value2 = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value2 = '1.309.29'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value2 = '  -5.13%  '
$ws.Range('B37').NumberFormat = '@'
$ws.Range('B37').Value2 = 'TrustWalletToken'
$ws.Range('C37').NumberFormat = '@'
$ws.Range('C37').Value2 = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value2 = '1.06'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value2 = '  -0.16%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value2 = '2.37'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value2 = '  +2.51%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value2 = '  +0.13%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value2 = '  +4.86%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value2 = '14.65'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value2 = '  +10.13%  '
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value2 = 'HuobiToken'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value2 = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value2 = '2.45'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value2 = '  +1.05%  '
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value2 = 'Aave'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value2 = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value2 = '82.53'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value2 = '  +0.58%  '
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value2 = 'MXToken'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value2 = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value2 = '2.80'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value2 = '  +0.44%  '
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value2 = 'ARBITRUM'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value2 = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value2 = '0.943'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value2 = '  -0.40%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value2 = '0.0519'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value2 = '  +3.70%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value2 = '1.969.55'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value2 = '  -0.01%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value2 = '5.73'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value2 = '  -3.86%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value2 = '  +0.07%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value2 = '102.00'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value2 = '  -0.97%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value2 = '  -0.45%  '
